# Fix Training Data Issue (#48)
#
# The "Date" column had been populated with the source file's name
# ("5-28-2007-08") instead of the actual game date. Because of the way the
# NBA stats site showed dates, the data had drifted one day off, so every
# row's Date value is corrected here to the real ISO date "2008-05-28".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$lastCol = $ws.UsedRange.Columns.Count

# Locate the "Date" column by its header text instead of hard-coding BF.
$dateCol = -1
for ($c = 1; $c -le $lastCol; $c++) {
    if ($ws.Cells.Item(1, $c).Value2 -eq "Date") {
        $dateCol = $c
        break
    }
}

$oldValue = "5-28-2007-08"
$newValue = "2008-05-28"

if ($dateCol -ne -1) {
    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, $dateCol)
        if ($cell.Value2 -eq $oldValue) {
            # Typing a date-shaped string normally makes Excel auto-convert
            # the cell to a date serial. Force text storage instead (so the
            # corrected value is kept as a literal string, matching the
            # original "Date" column's text format), then drop the cell
            # back to the workbook's default "Normal" style so no stray
            # number-format survives on the cell.
            $cell.NumberFormat = "@"
            $cell.Value = $newValue
            $cell.Style = "Normal"
        }
    }
}
